# Fill in the evidence values across the "A1".."A6" claim sheets, then
# move the active tab/selection state to match the saved workbook.

$wb = $excel.ActiveWorkbook

# --- A1 sheet: TxHash / ClassID -------------------------------------------------
$ws = $wb.Worksheets.Item("A1")
$ws.Cells.Item(2,1).Value = "C560BD65D1DF934CF8F531E56B761EC099A3843610FA3F4C7AD21BFB9A446F83`n"
$ws.Cells.Item(2,2).Value = "eknnft"
$ws.Cells.Item(2,1).WrapText = $true

# --- A2 sheet: TxHash / ClassID / NFTID (two NFTs minted) -----------------------
$ws = $wb.Worksheets.Item("A2")
$ws.Cells.Item(2,1).Value = "B18437988567B85CFE6F3A1D2F5D5E3A07E85A67811C2DB94CFB7A2343AC92EF"
$ws.Cells.Item(2,2).Value = "eknnft"
$ws.Cells.Item(2,3).Value = "nft00001"
$ws.Cells.Item(3,1).Value = "B120D412948B3CFF0D49BFE14F5C6615B7E3105BBEB8835802EAC5E81259F5A3`n"
$ws.Cells.Item(3,2).Value = "eknnft"
$ws.Cells.Item(3,3).Value = "nft00002"
$ws.Range("A3:C3").WrapText = $true
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- A3 sheet: TxHash / ClassID / NFTID / ChainID -------------------------------
$ws = $wb.Worksheets.Item("A3")
$ws.Cells.Item(2,1).Value = "C90832490D1965C38C31549B8E679012CFEF8C443808B43572756B9315E826DC"
$ws.Cells.Item(2,2).Value = "juno1zuqnqp0ytffafdt43sa34j5mk0yg8zvccgpalfulwlh57pf0x4rqglacue"
$ws.Cells.Item(2,3).Value = "nft00001"
$ws.Cells.Item(2,4).Value = "uni-6"

# --- A4 sheet: TxHash / ClassID / NFTID / ChainID -------------------------------
$ws = $wb.Worksheets.Item("A4")
$ws.Cells.Item(2,1).Value = "3124DB3F81484BBC37925843D94FC4D85270FF9016E54E7F80677DD3E2522D02"
$ws.Cells.Item(2,2).Value = "ibc/F2EF0D84F058C881E07A2311AD0D863004623E554D59805AFB931CB19CE1FC02`n"
$ws.Cells.Item(2,3).Value = "nft00002"
$ws.Cells.Item(2,4).Value = "gon-flixnet-1"
$ws.Cells.Item(2,2).WrapText = $true

# --- A5 sheet: TxHash / ClassID / NFTID / ChainID -------------------------------
$ws = $wb.Worksheets.Item("A5")
$ws.Cells.Item(2,1).Value = "5563956576FB21BA987431C2463F7E3C9506AF3C47655828AE6669EF4E7487E7"
$ws.Cells.Item(2,2).Value = "juno1zuqnqp0ytffafdt43sa34j5mk0yg8zvccgpalfulwlh57pf0x4rqglacue"
$ws.Cells.Item(2,3).Value = "nft00001"
$ws.Cells.Item(2,4).Value = "uni-6"

# --- A6 sheet: TxHash / ClassID / NFTID / ChainID -------------------------------
$ws = $wb.Worksheets.Item("A6")
$ws.Cells.Item(2,1).Value = "E011FA1A5AEC894F4F26FFC2527E506BE8212338540DD3579D8656C8CF9342B0"
$ws.Cells.Item(2,2).Value = "ibc/F2EF0D84F058C881E07A2311AD0D863004623E554D59805AFB931CB19CE1FC02"
$ws.Cells.Item(2,3).Value = "nft00002"
$ws.Cells.Item(2,4).Value = "gon-flixnet-1"

# --- Selections per sheet, matching the saved workbook state -------------------
$wb.Worksheets.Item("A1").Range("B2").Select()
$wb.Worksheets.Item("A2").Range("C3").Select()
$wb.Worksheets.Item("A3").Range("D2").Select()
$wb.Worksheets.Item("A4").Range("D2").Select()
$wb.Worksheets.Item("A5").Range("C3").Select()
$wb.Worksheets.Item("A6").Range("D2").Select()
$wb.Worksheets.Item("Info").Range("B2").Select()

# Move the active/selected tab to "A6" (was "Info").
$wb.Worksheets.Item("A6").Activate()
